{"js": "// Update the date heading and the 25 multiplication-table answer cells.\n// Cell edits are addressed by (row, col) rather than by searching for the\n// old text, because several new values collide with old values that live\n// elsewhere in the table (e.g. \"44\u00d789=3916\" is both an old value at\n// row 0 col 3 AND the new value for row 9 col 3) \u2014 a plain text search/\n// replace done sequentially would risk re-matching an already-updated cell.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// 1) Date heading paragraph: \"2026-02-04 Wednesday\" -> \"2026-02-05 Thursday\"\nconst firstPara = paragraphs.items[0];\nfirstPara.load(\"text\");\nawait context.sync();\nif (firstPara.text.trim() === \"2026-02-04 Wednesday\") {\n  firstPara.insertText(\"2026-02-05 Thursday\", \"Replace\");\n} else {\n  // Fallback: locate it anywhere in the body if layout ever shifts.\n  const hits = body.search(\"2026-02-04 Wednesday\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(\"2026-02-05 Thursday\", \"Replace\");\n  }\n}\n\n// 2) Table cell values, addressed by (row, col) \u2014 row indices 0, 4, 9, 14, 19\n//    are the five populated rows (the rows in between are blank spacers).\nconst table = tables.items[0];\nconst updates = [\n  [0, 0, \"47\u00d771=3337\"],\n  [0, 1, \"88\u00d775=6600\"],\n  [0, 2, \"96\u00d723=2208\"],\n  [0, 3, \"36\u00d771=2556\"],\n  [0, 4, \"30\u00d766=1980\"],\n  [4, 0, \"50\u00d765=3250\"],\n  [4, 1, \"58\u00d736=2088\"],\n  [4, 2, \"44\u00d773=3212\"],\n  [4, 3, \"47\u00d759=2773\"],\n  [4, 4, \"92\u00d753=4876\"],\n  [9, 0, \"68\u00d719=1292\"],\n  [9, 1, \"91\u00d726=2366\"],\n  [9, 2, \"77\u00d732=2464\"],\n  [9, 3, \"44\u00d789=3916\"],\n  [9, 4, \"19\u00d726=494\"],\n  [14, 0, \"68\u00d775=5100\"],\n  [14, 1, \"17\u00d734=578\"],\n  [14, 2, \"41\u00d794=3854\"],\n  [14, 3, \"85\u00d762=5270\"],\n  [14, 4, \"97\u00d727=2619\"],\n  [19, 0, \"86\u00d778=6708\"],\n  [19, 1, \"11\u00d782=902\"],\n  [19, 2, \"30\u00d754=1620\"],\n  [19, 3, \"54\u00d796=5184\"],\n  [19, 4, \"35\u00d779=2765\"],\n];\n\nfor (const [row, col, newText] of updates) {\n  table.getCell(row, col).value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 multiplication-table answer cells.\n# Cell edits are addressed by (row, col) rather than by a document-wide\n# Find/Replace, because several new values collide with old values that\n# live elsewhere in the table (e.g. \"44\u00d789=3916\" is both an old value at\n# row 1 col 4 AND the new value for row 10 col 4 in 1-based Word indexing)\n# \u2014 a sequential global replace would risk re-matching an already-updated\n# cell. Setting each cell's Range.Text directly keeps every edit scoped to\n# its own cell.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph: \"2026-02-04 Wednesday\" -> \"2026-02-05 Thursday\"\n$d.Paragraphs.Item(1).Range.Text = \"2026-02-05 Thursday\"\n\n# 2) Table cell values, addressed by 1-based (row, col). Rows 1, 5, 10, 15, 20\n#    are the five populated rows (the rows in between are blank spacers).\n$tbl = $d.Tables.Item(1)\n\n$tbl.Cell(1, 1).Range.Text = \"47\u00d771=3337\"\n$tbl.Cell(1, 2).Range.Text = \"88\u00d775=6600\"\n$tbl.Cell(1, 3).Range.Text = \"96\u00d723=2208\"\n$tbl.Cell(1, 4).Range.Text = \"36\u00d771=2556\"\n$tbl.Cell(1, 5).Range.Text = \"30\u00d766=1980\"\n\n$tbl.Cell(5, 1).Range.Text = \"50\u00d765=3250\"\n$tbl.Cell(5, 2).Range.Text = \"58\u00d736=2088\"\n$tbl.Cell(5, 3).Range.Text = \"44\u00d773=3212\"\n$tbl.Cell(5, 4).Range.Text = \"47\u00d759=2773\"\n$tbl.Cell(5, 5).Range.Text = \"92\u00d753=4876\"\n\n$tbl.Cell(10, 1).Range.Text = \"68\u00d719=1292\"\n$tbl.Cell(10, 2).Range.Text = \"91\u00d726=2366\"\n$tbl.Cell(10, 3).Range.Text = \"77\u00d732=2464\"\n$tbl.Cell(10, 4).Range.Text = \"44\u00d789=3916\"\n$tbl.Cell(10, 5).Range.Text = \"19\u00d726=494\"\n\n$tbl.Cell(15, 1).Range.Text = \"68\u00d775=5100\"\n$tbl.Cell(15, 2).Range.Text = \"17\u00d734=578\"\n$tbl.Cell(15, 3).Range.Text = \"41\u00d794=3854\"\n$tbl.Cell(15, 4).Range.Text = \"85\u00d762=5270\"\n$tbl.Cell(15, 5).Range.Text = \"97\u00d727=2619\"\n\n$tbl.Cell(20, 1).Range.Text = \"86\u00d778=6708\"\n$tbl.Cell(20, 2).Range.Text = \"11\u00d782=902\"\n$tbl.Cell(20, 3).Range.Text = \"30\u00d754=1620\"\n$tbl.Cell(20, 4).Range.Text = \"54\u00d796=5184\"\n$tbl.Cell(20, 5).Range.Text = \"35\u00d779=2765\"\n"}
